# ------------------------------------------------------------------
# 06186-中国飞鹤.xlsx  — add a "2022-Q4" quarter
#   1. Insert a new "2022-Q4" sheet (fund-holdings detail) right
#      after "总计", built by cloning the "2020-Q4" sheet (same
#      layout/styling) and replacing its data.
#   2. Insert a new summary row for "2022-Q4" at the top of the
#      "总计" data table, renumbering the existing index column.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "2022-Q4" worksheet -------------------------
# Clone "2020-Q4" (last sheet) so headers/borders/bold styling (s=2)
# match the other quarter sheets, then drop it in right after "总计".
$template = $wb.Worksheets.Item(9)
$anchor = $wb.Worksheets.Item(2)
$template.Copy($anchor, $null)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Template (2020-Q4) had 17 data+header rows, 2022-Q4 only needs 12.
$q4.Range("A13:H17").EntireRow.Delete()

# Header differs slightly for this sheet ("基金规模" not "基金金额").
$q4.Range("D1").Value = "基金规模"

# D:G hold decimal-looking values that must stay text (e.g. "3.80"
# must not collapse to 3.8); B holds zero-padded fund codes that must
# not become numbers either. Force text formatting before writing.
$q4.Range("B2:G12").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "012892"
$q4.Range("C2").Value = "安信优质企业三年持有混合A"
$q4.Range("D2").Value = "15.58"
$q4.Range("E2").Value = "90.59"
$q4.Range("F2").Value = "3.80"
$q4.Range("G2").Value = "0.5920"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "012893"
$q4.Range("C3").Value = "安信优质企业三年持有混合C"
$q4.Range("D3").Value = "7.05"
$q4.Range("E3").Value = "90.59"
$q4.Range("F3").Value = "3.80"
$q4.Range("G3").Value = "0.2679"
$q4.Range("H3").Value = 10

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "011856"
$q4.Range("C4").Value = "安信均衡成长18个月持有混合A"
$q4.Range("D4").Value = "4.90"
$q4.Range("E4").Value = "92.55"
$q4.Range("F4").Value = "4.59"
$q4.Range("G4").Value = "0.2249"
$q4.Range("H4").Value = 9

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "008227"
$q4.Range("C5").Value = "宝盈研究精选混合A"
$q4.Range("D5").Value = "3.93"
$q4.Range("E5").Value = "91.85"
$q4.Range("F5").Value = "4.35"
$q4.Range("G5").Value = "0.1710"
$q4.Range("H5").Value = 10

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "014506"
$q4.Range("C6").Value = "博时成长臻选混合A"
$q4.Range("D6").Value = "3.24"
$q4.Range("E6").Value = "87.77"
$q4.Range("F6").Value = "3.50"
$q4.Range("G6").Value = "0.1134"
$q4.Range("H6").Value = 9

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "009880"
$q4.Range("C7").Value = "安信成长动力一年持有期混合"
$q4.Range("D7").Value = "1.94"
$q4.Range("E7").Value = "93.22"
$q4.Range("F7").Value = "4.13"
$q4.Range("G7").Value = "0.0801"
$q4.Range("H7").Value = 9

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "008228"
$q4.Range("C8").Value = "宝盈研究精选混合C"
$q4.Range("D8").Value = "1.01"
$q4.Range("E8").Value = "91.85"
$q4.Range("F8").Value = "4.35"
$q4.Range("G8").Value = "0.0439"
$q4.Range("H8").Value = 10

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "011858"
$q4.Range("C9").Value = "安信消费升级一年持有期混合A"
$q4.Range("D9").Value = "0.84"
$q4.Range("E9").Value = "88.78"
$q4.Range("F9").Value = "3.99"
$q4.Range("G9").Value = "0.0335"
$q4.Range("H9").Value = 2

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "011857"
$q4.Range("C10").Value = "安信均衡成长18个月持有混合C"
$q4.Range("D10").Value = "0.26"
$q4.Range("E10").Value = "92.55"
$q4.Range("F10").Value = "4.59"
$q4.Range("G10").Value = "0.0119"
$q4.Range("H10").Value = 9

$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "011859"
$q4.Range("C11").Value = "安信消费升级一年持有期混合C"
$q4.Range("D11").Value = "0.13"
$q4.Range("E11").Value = "88.78"
$q4.Range("F11").Value = "3.99"
$q4.Range("G11").Value = "0.0052"
$q4.Range("H11").Value = 2

$q4.Range("A12").Value = 10
$q4.Range("B12").Value = "014507"
$q4.Range("C12").Value = "博时成长臻选混合C"
$q4.Range("D12").Value = "0.14"
$q4.Range("E12").Value = "87.77"
$q4.Range("F12").Value = "3.50"
$q4.Range("G12").Value = "0.0049"
$q4.Range("H12").Value = 9

# --- 2. Add the "2022-Q4" row to the "总计" summary sheet ----------
$summary = $wb.Worksheets.Item(1)

# Push the existing data rows (old row 2..9) down to row 3..10.
$summary.Rows.Item(2).Insert()

# The inserted row copies the header row's bold/bordered style —
# reset B2:D2 back to the plain "Normal" style used by data rows.
$summary.Range("B2:D2").Style = "Normal"

# A2 needs the same index-column style (s=2) as the rest of column A;
# copy formatting from A3 (the row right below, already shifted down).
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 11
$summary.Range("D2").Value = 1.55

# Renumber the index column (A) for every pre-existing row: it is a
# plain 0-based row counter, so each one shifts up by one.
for ($r = 3; $r -le 10; $r++) {
    $cell = $summary.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}
